$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Acre
$ws.Range("A2").Value = "Acre"
$ws.Range("B2").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C2").Value = 1.2

# Row 3: Maranhão
$ws.Range("A3").Value = "Maranhão"
$ws.Range("B3").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C3").Value = 0.8999999999999995

# Row 4: Distrito Federal
$ws.Range("A4").Value = "Distrito Federal"
$ws.Range("B4").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C4").Value = -0.09999999999999964

# Row 5: Mato Grosso
$ws.Range("A5").Value = "Mato Grosso"
$ws.Range("B5").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C5").Value = -0.1000000000000001

# Row 6: Rondônia
$ws.Range("A6").Value = "Rondônia"
$ws.Range("B6").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C6").Value = -0.1999999999999997

# Row 7: Rio Grande do Sul
$ws.Range("A7").Value = "Rio Grande do Sul"
$ws.Range("B7").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C7").Value = -0.3000000000000007

# Row 8: Sergipe (name unchanged)
$ws.Range("B8").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C8").Value = -1.4
$ws.Range("D8").Value = "20º"

# Row 9: Nordeste (name unchanged)
$ws.Range("B9").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C9").Value = -2.100000000000001

# Row 10: Brasil (name unchanged)
$ws.Range("B10").Value = "Diferença 2024/03 - 2023/03"
$ws.Range("C10").Value = -1.3
